$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Black Desert gear numbers were refreshed (AP/AAP/DP/GS for row 2).
# The diff shows D2 and G2 flipping from numeric storage to literal text
# (matching the existing text storage already used by E2/F2), with the
# new values. Prefixing with an apostrophe forces Excel to store the
# value as literal text instead of re-parsing it as a number.
$ws.Range("D2").Value = "'198"
$ws.Range("E2").Value = "'208"
$ws.Range("F2").Value = "'298"
$ws.Range("G2").Value = "'501"
